$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Refresh the "today" date fields shown in the footer of the slide master
#    and every slide layout (8/4/2025 -> 8/5/2025).
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "8/5/2025"
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2. Append two new "Title and Content" slides at the end of the deck.
# ---------------------------------------------------------------------------

# --- Slide 9: "Google fonts" ------------------------------------------------
$s9 = $p.Slides.Add(9, 2)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Google fonts"

$body9 = $s9.Shapes.Item(2).TextFrame.TextRange
$body9.Text = "<link "
$body9.InsertAfter("rel")
$body9.InsertAfter('="')
$body9.InsertAfter("preconnect")
$body9.InsertAfter('" ')
$body9.InsertAfter("href")
$body9.InsertAfter('="https://fonts.googleapis.com">')
$body9.InsertAfter("`r")
$body9.InsertAfter("<link ")
$body9.InsertAfter("rel")
$body9.InsertAfter('="')
$body9.InsertAfter("preconnect")
$body9.InsertAfter('" ')
$body9.InsertAfter("href")
$body9.InsertAfter('="https://fonts.gstatic.com" ')
$body9.InsertAfter("crossorigin")
$body9.InsertAfter(">")
$body9.InsertAfter("`r")
$body9.InsertAfter("<link ")
$body9.InsertAfter("href")
$body9.InsertAfter('="https://fonts.googleapis.com/css2?family=Open+Sans:ital,wght@0,300..800;1,300..800&family=Playwrite+HU:wght@100..400&family=Radio+Canada:ital,wght@0,300..700;1,300..700&display=swap" ')
$body9.InsertAfter("rel")
$body9.InsertAfter('="stylesheet">')

# --- Slide 10: "Color picker" -----------------------------------------------
$s10 = $p.Slides.Add(10, 2)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Color picker"

$body10 = $s10.Shapes.Item(2).TextFrame.TextRange
$body10.Text = "https://pickcoloronline.com/"
$body10.ActionSettings(1).Hyperlink.Address = "https://pickcoloronline.com/"
